$wb = $excel.ActiveWorkbook
$wsProperty = $wb.Worksheets.Item("Property")
$wsRecord = $wb.Worksheets.Item("Record")

# Rename the "View" header (shared string) to "Cache" everywhere it appears.
$wsProperty.Range("F1").Value = "Cache"
$wsRecord.Range("G1").Value = "Cache"

# Default the "Cache" (formerly "View") column to FALSE for all data rows.
$wsProperty.Range("F2:F26").Value = $false

# Rows 24-26 previously carried an extra (row-banding) fill on column F only
# because they sat outside the TRUE/FALSE validated block; now that they are
# treated like every other data row, clear that leftover formatting so F
# matches the plain look of F2:F23 (copy the unstyled format from F2).
$wsProperty.Range("F2").Copy() | Out-Null
$wsProperty.Range("F24:F26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Extend the TRUE/FALSE list validation so it also covers F24:F26 (it used
# to stop at F23, leaving a gap before the F27:F1048576 block).
$combined = $excel.Union($wsProperty.Range("F2:F26"), $wsProperty.Range("F27:F1048576"))
$combined.Validation.Delete() | Out-Null
$combined.Validation.Add(3, 1, 1, "TRUE,FALSE") | Out-Null

# Property sheet becomes the active/selected sheet, with F2:F26 selected
# and the view scrolled back to the top-left (A1).
$wsProperty.Range("A1").Select() | Out-Null
$wsProperty.Range("F2:F26").Select() | Out-Null
$wsProperty.Activate()
